# 自动更新价格数据: insert a new latest-date row at the top of the data
# (row 2), pushing the existing historical rows down by one. The new
# row carries the same commodity values as the previous latest row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new data row right below the header row.
$ws.Rows("2:2").Insert()

# Force column A to be treated as plain text so the date string isn't
# reinterpreted as a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-29"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Strip any formatting the row-insert/number-format step may have
# picked up so the new row matches the unstyled look of the other
# data rows.
$ws.Range("A2:D2").Style = "Normal"
